$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 119, shifting existing rows 119:238 down to 120:239.
$ws.Rows.Item(119).Insert()

# Populate the freshly inserted row 119 with the new weekly record.
$ws.Range("A119").Value = 3
$ws.Range("B119").Value = "Femacal de La Calera"
$ws.Range("C119").Value = "Coquimbo"
$ws.Range("D119").Value = 44586
$ws.Range("E119").Value = 5
$ws.Range("F119").Value = "Fruta"
$ws.Range("G119").Value = 100101
$ws.Range("H119").Value = "Berries"
$ws.Range("I119").Value = 100112025
$ws.Range("J119").Value = "Frutilla"
$ws.Range("K119").Value = "Sin especificar"
$ws.Range("L119").Value = "Primera"
$ws.Range("M119").Value = 120
$ws.Range("N119").Value = 5000
$ws.Range("O119").Value = 5000
$ws.Range("P119").Value = 5000
$ws.Range("Q119").Value = "$/bandeja 7 kilos"
$ws.Range("R119").Value = "Provincia de Melipilla"
$ws.Range("S119").Value = 714
$ws.Range("T119").Value = 7
